# Apply updated crypto price/volume data (GitHub Actions refresh).
# Leading "'" forces Excel to store the literal text verbatim
# (avoids "1.00" -> 1, "0.0370" -> 0.037 numeric coercion) without
# touching the cell NumberFormat/style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.680.46"
$ws.Range("E2").Value = "'  +0.27%  "

$ws.Range("D3").Value = "'3.137.17"
$ws.Range("E3").Value = "'  +0.60%  "

$ws.Range("E4").Value = "'  -0.16%  "

$ws.Range("D5").Value = "'586.45"
$ws.Range("E5").Value = "'  +0.01%  "

$ws.Range("D6").Value = "'145.41"
$ws.Range("E6").Value = "'  -0.02%  "

$ws.Range("E7").Value = "'  -0.07%  "

$ws.Range("D8").Value = "'3.134.56"
$ws.Range("E8").Value = "'  +0.75%  "

$ws.Range("E9").Value = "'  -0.55%  "

$ws.Range("E10").Value = "'  +5.90%  "

$ws.Range("D11").Value = "'5.72"
$ws.Range("E11").Value = "'  -1.15%  "

$ws.Range("E12").Value = "'  -2.55%  "

$ws.Range("E13").Value = "'  -0.63%  "

$ws.Range("D14").Value = "'36.95"
$ws.Range("E14").Value = "'  +3.99%  "

$ws.Range("B15").Value = "'TRON"
$ws.Range("C15").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "'0.121"
$ws.Range("E15").Value = "'  -1.70%  "

$ws.Range("B16").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "'3.656.98"
$ws.Range("E16").Value = "'  +0.64%  "

$ws.Range("B17").Value = "'WrappedBTC"
$ws.Range("C17").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'63.538.50"
$ws.Range("E17").Value = "'  +0.26%  "

$ws.Range("B18").Value = "'WrappedEther"
$ws.Range("C18").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.136.79"
$ws.Range("E18").Value = "'  +0.77%  "

$ws.Range("D20").Value = "'463.51"
$ws.Range("E20").Value = "'  -0.92%  "

$ws.Range("E21").Value = "'  +0.65%  "

$ws.Range("E22").Value = "'  +0.35%  "

$ws.Range("D23").Value = "'7.42"
$ws.Range("E23").Value = "'  -1.51%  "

$ws.Range("D24").Value = "'12.94"
$ws.Range("E24").Value = "'  -2.64%  "

$ws.Range("D25").Value = "'81.13"
$ws.Range("E25").Value = "'  -0.94%  "

$ws.Range("D26").Value = "'2.20"
$ws.Range("E26").Value = "'  +1.62%  "

$ws.Range("E27").Value = "'  -0.04%  "

$ws.Range("E28").Value = "'  +6.48%  "

$ws.Range("E29").Value = "'  +0.26%  "

$ws.Range("E30").Value = "'  -0.68%  "

$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "'  -0.19%  "

$ws.Range("D32").Value = "'6.97"
$ws.Range("E32").Value = "'  +1.92%  "

$ws.Range("E33").Value = "'  -0.28%  "

$ws.Range("E34").Value = "'  +0.41%  "

$ws.Range("E35").Value = "'  -3.22%  "

$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = "'  -0.37%  "

$ws.Range("E37").Value = "'  -5.31%  "

$ws.Range("E38").Value = "'  -0.59%  "

$ws.Range("E39").Value = "'  -1.57%  "

$ws.Range("E40").Value = "'  +0.50%  "

$ws.Range("D41").Value = "'439.71"
$ws.Range("E41").Value = "'  -0.96%  "

$ws.Range("D42").Value = "'8.82"
$ws.Range("E42").Value = "'  +1.28%  "

$ws.Range("D43").Value = "'0.0370"
$ws.Range("E43").Value = "'  +0.30%  "

$ws.Range("D44").Value = "'2.902.66"
$ws.Range("E44").Value = "'  -0.30%  "

$ws.Range("D45").Value = "'0.277"
$ws.Range("E45").Value = "'  -0.92%  "

$ws.Range("E46").Value = "'  -2.59%  "

$ws.Range("D47").Value = "'37.02"
$ws.Range("E47").Value = "'  +2.91%  "

$ws.Range("D48").Value = "'125.63"
$ws.Range("E48").Value = "'  +1.55%  "

$ws.Range("E50").Value = "'  -1.11%  "

$ws.Range("D51").Value = "'24.26"
$ws.Range("E51").Value = "'  -1.64%  "
